# Rename the "station" column header to "name" in the tram station list.
# Renaming the header cell of an Excel Table (ListObject) column updates
# both the worksheet cell and the table's column definition in one go.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"
